$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay stored as text,
# matching the source data (inline strings like "64.106.47", "571.90").
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "64.146.93"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3
$ws.Range("D3").Value = "3.423.58"
$ws.Range("E3").Value = "  +0.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "571.90"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6
$ws.Range("D6").Value = "161.13"
$ws.Range("E6").Value = "  +1.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "3.426.21"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  -5.84%  "

# Row 10
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +1.49%  "

# Row 11
$ws.Range("E11").Value = "  -1.79%  "

# Row 12
$ws.Range("E12").Value = "  -2.88%  "

# Row 13
$ws.Range("D13").Value = "4.012.67"
$ws.Range("E13").Value = "  +0.03%  "

# Row 14
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").Value = "27.01"
$ws.Range("E15").Value = "  -2.34%  "

# Row 16
$ws.Range("E16").Value = "  -6.65%  "

# Row 17
$ws.Range("D17").Value = "64.191.87"

# Row 18
$ws.Range("D18").Value = "3.425.12"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -4.08%  "

# Row 20
$ws.Range("D20").Value = "13.59"
$ws.Range("E20").Value = "  -1.82%  "

# Row 21
$ws.Range("D21").Value = "376.89"
$ws.Range("E21").Value = "  -0.70%  "

# Row 22
$ws.Range("D22").Value = "7.81"
$ws.Range("E22").Value = "  -2.35%  "

# Row 23
$ws.Range("E23").Value = "  -0.49%  "

# Row 24
$ws.Range("D24").Value = "71.27"
$ws.Range("E24").Value = "  -0.80%  "

# Row 25
$ws.Range("D25").Value = "0.519"
$ws.Range("E25").Value = "  -5.12%  "

# Row 26
$ws.Range("D26").Value = "0.0000116"
$ws.Range("E26").Value = "  -2.17%  "

# Row 27
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  -3.92%  "

# Row 28
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("D30").Value = "6.01"
$ws.Range("E30").Value = "  -2.39%  "

# Row 31
$ws.Range("E31").Value = "  -3.63%  "

# Row 33
$ws.Range("D33").Value = "22.89"
$ws.Range("E33").Value = "  -1.28%  "

# Row 34
$ws.Range("D34").Value = "7.10"
$ws.Range("E34").Value = "  +1.40%  "

# Row 35
$ws.Range("E35").Value = "  -3.94%  "

# Row 36
$ws.Range("D36").Value = "159.68"
$ws.Range("E36").Value = "  -0.70%  "

# Row 37
$ws.Range("D37").Value = "0.857"
$ws.Range("E37").Value = "  +11.18%  "

# Row 38
$ws.Range("E38").Value = "  -4.40%  "

# Row 39
$ws.Range("D39").Value = "2.804.95"
$ws.Range("E39").Value = "  -2.96%  "

# Row 40
$ws.Range("D40").Value = "0.0726"
$ws.Range("E40").Value = "  -3.30%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "25.84"
$ws.Range("E41").Value = "  -1.49%  "

# Row 42
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "42.91"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("D43").Value = "6.44"
$ws.Range("E43").Value = "  -4.14%  "

# Row 44
$ws.Range("D44").Value = "4.42"
$ws.Range("E44").Value = "  -2.42%  "

# Row 45
$ws.Range("D45").Value = "25.99"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46
$ws.Range("D46").Value = "0.0306"
$ws.Range("E46").Value = "  -3.03%  "

# Row 47
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +7.41%  "

# Row 48
$ws.Range("D48").Value = "335.14"
$ws.Range("E48").Value = "  +5.36%  "

# Row 50
$ws.Range("D50").Value = "6.31"
$ws.Range("E50").Value = "  -2.96%  "

# Row 51
$ws.Range("E51").Value = "  -3.00%  "

# Restore default styling on the Price column (value stays text; style index reverts).
$ws.Range("D2:D51").Style = "Normal"
